$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New worker data for rows 16-36 (replaces the old two-worker table with the
# updated "ANGIE OCHOA NAVAS" records - part 1 of the new estado de cuenta).
$rows = @(
    @{Row=16; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2305"; F=52753; G=1318840},
    @{Row=17; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2306"; F=52753; G=1318840},
    @{Row=18; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2307"; F=52753; G=1318840},
    @{Row=19; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2308"; F=52753; G=1318840},
    @{Row=20; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2309"; F=52753; G=1318840},
    @{Row=21; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2310"; F=52753; G=1318840},
    @{Row=22; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2311"; F=52753; G=1318840},
    @{Row=23; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2312"; F=48720; G=1218000},
    @{Row=24; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2312"; F=52753; G=1318840},
    @{Row=25; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2401"; F=48720; G=1218000},
    @{Row=26; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2401"; F=52753; G=1318840},
    @{Row=27; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2402"; F=48720; G=1218000},
    @{Row=28; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2402"; F=52753; G=1318840},
    @{Row=29; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2403"; F=48720; G=1218000},
    @{Row=30; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2403"; F=52753; G=1318840},
    @{Row=31; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2404"; F=48720; G=1218000},
    @{Row=32; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2404"; F=52753; G=1318840},
    @{Row=33; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2405"; F=48720; G=1218000},
    @{Row=34; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2405"; F=52753; G=1318840},
    @{Row=35; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2406"; F=19488; G=1218000},
    @{Row=36; B="CC"; C="1128054473"; D="ANGIE OCHOA NAVAS"; E="2406"; F=21101; G=1318840}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value2 = $r.B
    $ws.Cells.Item($row, 3).Value2 = $r.C
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 5).Value2 = $r.E
    $ws.Cells.Item($row, 6).Value2 = $r.F
    $ws.Cells.Item($row, 7).Value2 = $r.G
}
